# Fruta / hortaliza, semanal
# Inserts 4 new weekly price rows into the "Mandarina" dataset, pushing the
# existing data rows (225-265) down to (229-269).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 4 new rows by inserting above the current row 225.
$ws.Rows("225:228").Insert()

# Columns A,B,C,E,F,G,H,I,J are constant across every data row in this sheet.
$constA = 9
$constB = "Vega Central Mapocho de Santiago"
$constC = "Metropolitana"
$constE = 13
$constF = "Fruta"
$constG = 100102
$constH = "Cítricos"
$constI = 100102004
$constJ = "Mandarina"

# Data for the 4 newly inserted rows (225-228).
$newRows = @(
    @{ Row=225; D=44449; K="Murcott"; L="Especial"; M=450; N=9000; O=9000; P=9000; Q="$/caja 18 kilos";    R="Provincia de San Felipe de Aconcagua"; S=500; T=18 },
    @{ Row=226; D=44449; K="Murcott"; L="Primera";  M=380; N=5000; O=5000; P=5000; Q="$/bandeja 10 kilos"; R="Provincia de Limarí";                   S=500; T=10 },
    @{ Row=227; D=44449; K="Murcott"; L="Primera";  M=400; N=7000; O=7000; P=7000; Q="$/caja 18 kilos";    R="Provincia de San Felipe de Aconcagua"; S=389; T=18 },
    @{ Row=228; D=44449; K="Murcott"; L="Segunda";  M=380; N=5500; O=5500; P=5500; Q="$/caja 18 kilos";    R="Provincia de San Felipe de Aconcagua"; S=306; T=18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $constA
    $ws.Range("B$row").Value = $constB
    $ws.Range("C$row").Value = $constC
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $constE
    $ws.Range("F$row").Value = $constF
    $ws.Range("G$row").Value = $constG
    $ws.Range("H$row").Value = $constH
    $ws.Range("I$row").Value = $constI
    $ws.Range("J$row").Value = $constJ
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
}
